$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 17 (question 11 - "Mi a fo celja a tobbciklusu kepzes rendszernak?")
# This shifts rows 18-22 up to become 17-21, and prunes the now-unused
# shared strings automatically.
$ws.Rows(17).Delete()

# Renumber column B (the "sorszam" column) for the shifted rows so the
# sequence stays 11,12,13,14,15 and give them the new number-format style
# (General instead of Text) used for this id column going forward.
$ws.Range("B17").Value = 11
$ws.Range("B18").Value = 12
$ws.Range("B19").Value = 13
$ws.Range("B20").Value = 14
$ws.Range("B21").Value = 15

$idRange = $ws.Range("B17:B21")
$idRange.NumberFormat = "General"

# Update the selection left from editing (selection moved down while the
# view now starts scrolled back to the top of the sheet).
$ws.Range("A1").Select()
$ws.Range("I9").Select()
